$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.095.00"
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").Value = "'1.842.69"
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  -0.89%  '
$ws.Range("D4").Value = "'0.9978"
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = "'245.57"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +1.53%  '
$ws.Range("D6").Value = "'0.6975"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("D7").Value = "'0.9988"
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = "'0.07715"
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  -1.16%  '
$ws.Range("D9").Value = "'0.3058"
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  -1.55%  '
$ws.Range("D10").Value = "'23.56"
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  -1.47%  '
$ws.Range("D11").Value = "'0.07830"
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  +0.35%  '
$ws.Range("D12").Value = "'93.08"
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("D13").Value = "'1.838.32"
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  -1.25%  '
$ws.Range("D14").Value = "'5.122"
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("D15").Value = "'0.6852"
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  -0.83%  '
$ws.Range("D16").Value = "'6.617"
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("D17").Value = "'0.000008290"
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  -1.79%  '
$ws.Range("D18").Value = "'29.060.44"
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  -0.56%  '
$ws.Range("D19").Value = "'242.27"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  -3.08%  '
$ws.Range("D20").Value = "'2.077.76"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -1.79%  '
$ws.Range("D21").Value = "'12.75"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  -1.30%  '
$ws.Range("D23").Value = "'7.490"
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -1.39%  '
$ws.Range("D24").Value = "'0.9986"
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").Value = "'0.1510"
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  -1.41%  '
$ws.Range("D26").Value = "'159.09"
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  -1.12%  '
$ws.Range("D27").Value = "'8.814"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -0.94%  '
$ws.Range("D28").Value = "'18.24"
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  -1.68%  '
$ws.Range("D29").Value = "'1.540"
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -1.95%  '
$ws.Range("D30").Value = "'4.229"
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  -1.12%  '
$ws.Range("D31").Value = "'4.175"
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -1.79%  '
$ws.Range("D32").Value = "'1.199"
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  -1.21%  '
$ws.Range("D33").Value = "'0.05120"
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -2.12%  '
$ws.Range("D34").Value = "'0.7853"
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  +3.55%  '
$ws.Range("D35").Value = "'1.863"
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  -0.74%  '
$ws.Range("D36").Value = "'1.146"
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  -2.58%  '
$ws.Range("D37").Value = "'2.696"
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  -0.48%  '
$ws.Range("D38").Value = "'1.306.57"
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  +7.04%  '
$ws.Range("D39").Value = "'0.01865"
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  +0.14%  '
$ws.Range("E40").Value = '  -0.68%  '
$ws.Range("D41").Value = "'0.9480"
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  +5.33%  '
$ws.Range("D42").Value = "'6.141"
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  +5.63%  '
$ws.Range("D43").Value = "'107.82"
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -2.40%  '
$ws.Range("D44").Value = "'0.9989"
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").Value = "'9.706"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  +1.87%  '
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").Value = "'1.978.40"
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  -1.63%  '
$ws.Range("D48").Value = "'64.28"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -4.28%  '
$ws.Range("D49").Value = "'1.763"
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("D50").Value = "'0.00000000119"
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  -1.96%  '
$ws.Range("D51").Value = "'6.996"
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  -0.71%  '
